$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a "plain" number-looking string (e.g. "522.35")
# must be forced to Text so Excel does not re-interpret them as numeric values -
# matching the workbook author convention where every Price cell is inline text.
$textForceCells = @(
    'D5', 'D6', 'D16', 'D19', 'D21', 'D24', 'D26', 'D34', 'D37', 'D39', 'D43', 'D44', 'D47', 'D48'
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '59.017.29'
$ws.Range('E2').Value = '  +1.34%  '
$ws.Range('D3').Value = '2.587.80'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '522.35'
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('D6').Value = '139.46'
$ws.Range('E6').Value = '  -2.43%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.25%  '
$ws.Range('D9').Value = '2.601.02'
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('E10').Value = '  -4.41%  '
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('E13').Value = '  +3.56%  '
$ws.Range('D14').Value = '3.047.30'
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('D15').Value = '58.986.02'
$ws.Range('E15').Value = '  +1.46%  '
$ws.Range('D16').Value = '20.53'
$ws.Range('E16').Value = '  +0.90%  '
$ws.Range('D17').Value = '2.611.80'
$ws.Range('E17').Value = '  +0.66%  '
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('D19').Value = '338.14'
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').Value = '10.08'
$ws.Range('E21').Value = '  -1.90%  '
$ws.Range('E22').Value = '  +2.21%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = '66.23'
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('E25').Value = '  +1.61%  '
$ws.Range('D26').Value = '0.403'
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('E28').Value = '  +0.54%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '0.0₃0724'
$ws.Range('E30').Value = '  -2.72%  '
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('D34').Value = '149.08'
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('E35').Value = '  -1.37%  '
$ws.Range('E36').Value = '  -2.14%  '
$ws.Range('D37').Value = '36.79'
$ws.Range('E37').Value = '  +2.37%  '
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('D39').Value = '0.824'
$ws.Range('E39').Value = '  -1.53%  '
$ws.Range('E40').Value = '  -6.40%  '
$ws.Range('E41').Value = '  -0.73%  '
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = '272.03'
$ws.Range('E43').Value = '  -0.73%  '
$ws.Range('D44').Value = '10.78'
$ws.Range('E44').Value = '  +1.25%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D47').Value = '0.0517'
$ws.Range('E47').Value = '  -1.65%  '
$ws.Range('D48').Value = '18.42'
$ws.Range('E48').Value = '  -2.39%  '
$ws.Range('D49').Value = '1.963.06'
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('E51').Value = '  -0.39%  '
